$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new data rows right after the header/above the old row 112
# (shifts the existing rows 112:150 down to 114:152).
$ws.Rows.Item(112).Insert()
$ws.Rows.Item(113).Insert()

# New row 112: Black Amber / Primera, week of 2022-01-11 (serial 44572)
$ws.Cells.Item(112, 1).Value = 8
$ws.Cells.Item(112, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(112, 3).Value = "Coquimbo"
$ws.Cells.Item(112, 4).Value = 44572
$ws.Cells.Item(112, 5).Value = 4
$ws.Cells.Item(112, 6).Value = "Fruta"
$ws.Cells.Item(112, 7).Value = 100103
$ws.Cells.Item(112, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(112, 9).Value = 100103002
$ws.Cells.Item(112, 10).Value = "Ciruela"
$ws.Cells.Item(112, 11).Value = "Black Amber"
$ws.Cells.Item(112, 12).Value = "Primera"
$ws.Cells.Item(112, 13).Value = 20
$ws.Cells.Item(112, 14).Value = 290000
$ws.Cells.Item(112, 15).Value = 300000
$ws.Cells.Item(112, 16).Value = 295000
$ws.Cells.Item(112, 17).Value = "$/bins (450 kilos)"
$ws.Cells.Item(112, 18).Value = "Región Metropolitana"
$ws.Cells.Item(112, 19).Value = 656
$ws.Cells.Item(112, 20).Value = 450

# New row 113: Black Amber / Segunda, week of 2022-01-11 (serial 44572)
$ws.Cells.Item(113, 1).Value = 8
$ws.Cells.Item(113, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(113, 3).Value = "Coquimbo"
$ws.Cells.Item(113, 4).Value = 44572
$ws.Cells.Item(113, 5).Value = 4
$ws.Cells.Item(113, 6).Value = "Fruta"
$ws.Cells.Item(113, 7).Value = 100103
$ws.Cells.Item(113, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(113, 9).Value = 100103002
$ws.Cells.Item(113, 10).Value = "Ciruela"
$ws.Cells.Item(113, 11).Value = "Black Amber"
$ws.Cells.Item(113, 12).Value = "Segunda"
$ws.Cells.Item(113, 13).Value = 16
$ws.Cells.Item(113, 14).Value = 240000
$ws.Cells.Item(113, 15).Value = 250000
$ws.Cells.Item(113, 16).Value = 245000
$ws.Cells.Item(113, 17).Value = "$/bins (450 kilos)"
$ws.Cells.Item(113, 18).Value = "Región Metropolitana"
$ws.Cells.Item(113, 19).Value = 544
$ws.Cells.Item(113, 20).Value = 450

# Make sure the date cells keep the same date/time number format used elsewhere in column D
$ws.Cells.Item(112, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(113, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
